$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-6 (header row 1 stays unchanged).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @{ Row=2; A="ECs";               E=2; F=0.6666666666666666;  G=0.3893783333333333;  H=1.168135;          I=0.2294878876126519;  J=0.2521981668521232;  M=0.0237255; N=0.047451; Q=0.009238195647499999; R=0.055429173885;      S=0.2294878876126519;  T=0.2521981668521232 },
    @{ Row=3; A="FAPs";              E=3; F=1;                   G=0.7438396666666667;  H=2.231519;          I=0.4383967447919097;  J=0.4817807882613594;  M=0.0237255; N=0.047451; Q=0.0176479680115;       R=0.105887808069;      S=0.4383967447919097;  T=0.4817807882613594 },
    @{ Row=4; A="Inflammatory-Mac";  E=1; F=0.3333333333333333;  G=0.08827633333333333; H=0.264829;          I=0.05202741788283973; J=0.05717608694995093; M=0.0237255; N=0.047451; Q=0.0020944001465;       R=0.012566400879;      S=0.05202741788283973; T=0.05717608694995093 },
    @{ Row=5; A="MuSCs";             E=2; F=1;                   G=0.4583675;           H=0.9167350000000001; I=0.2701480291026951;  J=0.1979213759447163;  M=0.0237255; N=0.047451; Q=0.01087499812125;      R=0.04349999248500001; S=0.2701480291026951;  T=0.1979213759447163 },
    @{ Row=6; A="Resolving-Mac";     E=1; F=0.3333333333333333;  G=0.01686533333333333; H=0.050596;          I=0.009939920609903594; J=0.01092358199185028; M=0.0237255; N=0.047451; Q=0.000400138466;        R=0.002400830796;      S=0.009939920609903594; T=0.01092358199185028 }
)

foreach ($row in $data) {
    $r = $row.Row

    $ws.Cells.Item($r, 1).Value = $row.A          # A: Sending cluster
    $ws.Cells.Item($r, 2).Value = "Efna4"          # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value = "Epha5"          # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value = "MuSCs"          # D: Target cluster

    $ws.Cells.Item($r, 5).Value = $row.E           # E
    $ws.Cells.Item($r, 6).Value = $row.F           # F
    $ws.Cells.Item($r, 7).Value = $row.G           # G
    $ws.Cells.Item($r, 8).Value = $row.H           # H
    $ws.Cells.Item($r, 9).Value = $row.I           # I
    $ws.Cells.Item($r, 10).Value = $row.J          # J
    $ws.Cells.Item($r, 11).Value = 2               # K
    $ws.Cells.Item($r, 12).Value = 1               # L
    $ws.Cells.Item($r, 13).Value = $row.M          # M
    $ws.Cells.Item($r, 14).Value = $row.N          # N
    $ws.Cells.Item($r, 15).Value = 1               # O
    $ws.Cells.Item($r, 16).Value = 1               # P
    $ws.Cells.Item($r, 17).Value = $row.Q          # Q
    $ws.Cells.Item($r, 18).Value = $row.R          # R
    $ws.Cells.Item($r, 19).Value = $row.S          # S
    $ws.Cells.Item($r, 20).Value = $row.T          # T
}
